$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.580.24'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.743.25'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.64'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4927'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2676'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06283'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.742.94'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07041'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.73'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6148'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.584'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '78.09'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.12%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.592.04'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007288'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.55'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.966.68'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.579'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.719'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.261'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.54'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.46'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.427'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.763'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.043'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08039'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.737'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04629'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.40%  '
$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9996'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.612'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.24%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.016'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.53%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6385'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.064'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.53%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8973'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.46%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.425'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.003'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01504'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.88'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.05%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.428'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.79%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3917'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.871'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.89%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1182'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.02%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05399'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.45%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.54'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.49%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.794'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.264'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.33%  '
